# Apply the "add the root part" edit:
# - Remove the section_id column and the time/classroom_no/lesson/limit/day/instructor_id columns
# - Keep course_id, title, credits, dept_name
# - Add two data rows with course info
# - Update the selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "section_id" column (E) - columns shift left
$ws.Range("E1").EntireColumn.Delete() | Out-Null

# After that delete, the remaining unwanted columns (time, classroom_no, lesson,
# limit, day, instructor_id) now occupy E:J - remove them too
$ws.Range("E1:J1").EntireColumn.Delete() | Out-Null

# Fill in the new data rows
$ws.Range("A2").Value = "CCCC120001"
$ws.Range("B2").Value = "我不做人了"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "软件学院"

$ws.Range("A3").Value = "DDDD111111"
$ws.Range("B3").Value = "我死了"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "计算机学院"

# Move the active selection cell to match the edited file
$ws.Range("F4").Select() | Out-Null
